$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 2.228397109637203

$ws.Range("C4").ClearContents()

$ws.Range("C5").Value = -0.4626567965509643
$ws.Range("E5").Value = -0.2018858887078645

$ws.Range("C6").Value = -0.8792832172735632
$ws.Range("E6").Value = -0.1037449741509211

$ws.Range("E7").Value = 0.2108047537406454

$ws.Range("C8").Value = 0.9337833426867226

$ws.Range("E9").Value = 0.96842791562195

$ws.Range("C10").Value = 2.791140000794257
$ws.Range("E10").Value = 1.722110645261954

$ws.Range("E11").Value = 1.698182372097512

$ws.Range("C12").Value = 0.4451370000809529
$ws.Range("E12").Value = 0.6480763427742176

$ws.Range("C13").Value = -0.2674335569108788

$ws.Range("E14").Value = 1.223618887196509

$ws.Range("C15").Value = 2.038609866767915
$ws.Range("E15").Value = 1.325176859452348

$ws.Range("E17").Value = 1.522808462763692

$ws.Range("E18").Value = 1.582979977679555

$ws.Range("C19").Value = 2.246337373618967
$ws.Range("E19").Value = 1.693557061600948

$ws.Range("C20").Value = 2.2044495746113
$ws.Range("E20").Value = 1.401113624217065

$ws.Range("C21").Value = 1.777150434343544
$ws.Range("E21").Value = 1.905564797014669

$ws.Range("E22").Value = 2.357704431248386

$ws.Range("C23").Value = 1.953801996162019
$ws.Range("E23").Value = 2.2850393561338

$ws.Range("C24").Value = 1.741137453897323
$ws.Range("E24").Value = 2.08247707460909

$ws.Range("E25").Value = 2.263031582094888

$ws.Range("E26").Value = 2.36261304543155

$ws.Range("C27").Value = 1.831523226563148
$ws.Range("E27").Value = 2.069147625228918

$ws.Range("C28").Value = 1.562095320687429
$ws.Range("E28").Value = 1.845103901518907

$ws.Range("E29").Value = 1.670042545405948

$ws.Range("E30").Value = 0.7756897792100093

$ws.Range("C31").Value = 0.9716673922242069
$ws.Range("E31").Value = 1.221567067111828

$ws.Range("E32").Value = 0.02570757229445331

$ws.Range("C33").Value = -4.58200588280312
$ws.Range("E33").Value = -5.994170454492931

$ws.Range("C34").Value = -3.258619210312896
$ws.Range("E34").Value = -1.049961713694159

$ws.Range("C35").Value = -1.267247591471976

$ws.Range("C37").Value = -0.4671761491813142

$ws.Range("C38").Value = 0.4255262881966759
$ws.Range("E38").Value = 2.734996705911397

$ws.Range("C39").Value = 1.709116405894862
$ws.Range("E39").Value = -0.3788276656591538

$ws.Range("E42").Value = 1.813346177122321

$ws.Range("C43").Value = 1.845950461732082

$ws.Range("C44").Value = -0.2883789941992232

$ws.Range("C45").Value = -0.1703190263132703
$ws.Range("E45").Value = 1.429443202289193

$ws.Range("C46").Value = -0.2814561130375703

$ws.Range("E47").Value = 0.8057832540545151

$ws.Range("C49").Value = -0.9450815092640896
$ws.Range("E49").Value = 0.1338934428748884

$ws.Range("C50").Value = -0.6470065423293758
$ws.Range("E50").Value = 0.9227184786156251

$ws.Range("C51").Value = 2.137819896031878
$ws.Range("E51").Value = 0.3874019458836253

$ws.Range("C52").Value = 1.431852292002245

$ws.Range("C53").Value = 2.034296706251948
$ws.Range("E53").Value = 0.6096213037243281
